$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: remove "? Mini Project 03" from G17 ---
$ws.Range("G17").ClearContents()

# --- Row 19: remove "? Problem Set 02" from G19 ---
$ws.Range("G19").ClearContents()

# --- Row 20: add Mini Project 03 link to G20 ---
$ws.Range("G20").Value2 = "[Mini Project 03: Tudor network](https://docs.google.com/document/d/1XxtW5ZlPzRwHcbTHapKzeSZySq0axOVMnOGFtAnprEY/edit?usp=sharing)"

# --- Insert a new row at 23 for the Nov 11 holiday, shifting rows 23-31 down to 24-32 ---
$ws.Rows("23:23").Insert()

# New row 23 only has a Topic (column D) entry
$ws.Range("D23").Value2 = "NO CLASS (HOLIDAY)"

# After the insert, the rows read (old content shifted down by one):
#  row24: A24=[Social influence](#sec:socialinfluence)      D24=Social influence, herding, and cascades
#  row25: D25=Threshold models and complex contagion
#  row26: A26=[Dynamics: Complex contagion...]               D26=Complex contagion on networks
#  row27: D27=Complex contagion on networks, cont. + ...     G27=? Problem Set 3
#  row28: D28=NO CLASS
#  row29: D29=THANKSGIVING (NO CLASS)
#  row30: A30=[Cooperation](#sec:cooperation)                D30=Cooperation and networks
#  row31: D31=Wrap up
#  row32: D32=READING WEEK
#
# Target layout:
#  row26: ... G26=? Problem Set 3 (moved up from row27)
#  row27: D27=NO CLASS (moved up from row28)
#  row28: D28=THANKSGIVING (NO CLASS) (moved up from row29)
#  row29: D29=Complex contagion on networks, cont. + ... (moved down from row27, G cleared)
#  row30: A30 (from row30), D30=Cooperation and networks / Wrap Up (merge row30 + row31)
#  row31: D31=READING WEEK (from row32)

$contText = $ws.Range("D27").Value2

$ws.Range("G26").Value2 = $ws.Range("G27").Value2
$ws.Range("G27").ClearContents()

$ws.Range("D27").Value2 = $ws.Range("D28").Value2
$ws.Range("D28").Value2 = $ws.Range("D29").Value2
$ws.Range("D29").Value2 = $contText

$ws.Range("D30").Value2 = $ws.Range("D30").Value2 + " / Wrap Up"
$ws.Range("D31").Value2 = $ws.Range("D32").Value2

$ws.Rows("32:32").Delete()
